$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B52 becomes a real number instead of text "3"
$ws.Range("B52").Value = 3

# Add new row 53 with Ruilin's additional annotation
$ws.Range("A53").Value = "Ruilin"

# Keep B53 as text "2" (not auto-converted to a number)
$ws.Range("B53").Value = "'2"
$ws.Range("B53").Style = "Normal"

$ws.Range("C53").Value = "And I am still not convinced by the quality of the paper."
$ws.Range("D53").Value = "CRT"
$ws.Range("E53").Value = "OTH"
$ws.Range("F53").Value = "8e71ac7d-5b58-47e3-b02f-0c25e20406c1"
$ws.Range("G53").Value = "rJTGkKxAZ_annotated.xlsx"
$ws.Range("H53").Value = "And I am still not convinced by the quality of the paper."
